# Rename the element/Sector labels so they are aligned with the Baseline
# workbook's naming convention. These labels live in the shared-string
# table and are used as the column headers (row 3, columns D:G) on every
# yearly worksheet (2000 .. 2100).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Range("D3").Value2 -eq "Nd") {
        $ws.Range("D3").Value = "Neodymium"
    }
    if ($ws.Range("E3").Value2 -eq "Dy") {
        $ws.Range("E3").Value = "Dysprosium"
    }
    if ($ws.Range("F3").Value2 -eq "Cu") {
        $ws.Range("F3").Value = "Copper ores and concentrates"
    }
    if ($ws.Range("G3").Value2 -eq "Si") {
        $ws.Range("G3").Value = "Raw silicon"
    }
}

# A handful of cached "Raw silicon" (G) totals on specific yearly sheets
# carry a one-ULP rounding refresh so they stay in sync with the Baseline
# workbook after the relabeling.
$g7Updates = @{
    "2009" = -7419.091319003097
    "2020" = -6023034.283453048
    "2025" = -72575285.87201165
    "2026" = -82597841.80969585
    "2030" = -146271578.64677
    "2041" = -426400596.5710188
    "2042" = -500112075.5291855
    "2045" = -937876045.8612229
    "2047" = -1441021538.095722
    "2051" = -2539131669.523302
    "2097" = -273427420.4801398
}

foreach ($sheetName in $g7Updates.Keys) {
    $sheetNameText = [string]$sheetName
    $wb.Worksheets.Item($sheetNameText).Range("G7").Value = $g7Updates[$sheetName]
}
